$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new holiday-table row (date + ISO weekday number), continuing
# the existing day-by-day series through 2023-12-09.
$ws.Range("A70").Value = 20231209
$ws.Range("B70").Value = 6

# Match the author's final on-screen selection after adding the row.
$ws.Range("G70").Select()
